$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 75
$ws.Range("H75").Value = 40311.332
$ws.Range("J75").Value = 40311.332
$ws.Range("L75").Value = 40311.332
$ws.Range("N75").Value = -42183.332
# Row 78
$ws.Range("H78").Value = 40311.332
$ws.Range("J78").Value = 40311.332
$ws.Range("L78").Value = 120933.996
$ws.Range("N78").Value = -130293.996
# Row 137
$ws.Range("H137").Value = 1235.2858
$ws.Range("I137").Value = 903.8095
$ws.Range("K137").Value = 2711.4285
$ws.Range("M137").Value = -161.4285

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3038.6365
$ws.Range("I32").Value = 2479.84
$ws.Range("J32").Value = 6262.4614
$ws.Range("K32").Value = 2479.84
$ws.Range("L32").Value = 6262.4614
$ws.Range("M32").Value = -2192.84
$ws.Range("N32").Value = -6836.4614

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 88
$ws.Range("H88").Value = 15343
$ws.Range("J88").Value = 15343
$ws.Range("L88").Value = 15343
$ws.Range("N88").Value = -16155
# Row 91
$ws.Range("H91").Value = 15343
$ws.Range("J91").Value = 15343
$ws.Range("L91").Value = 15343
$ws.Range("N91").Value = -18151
# Row 105
$ws.Range("H105").Value = 2458.2917
$ws.Range("I105").Value = 2380.9048
$ws.Range("K105").Value = 2380.9048
$ws.Range("M105").Value = -633.9047999999998
# Row 126
$ws.Range("H126").Value = 39990.477
$ws.Range("J126").Value = 39990.477
$ws.Range("L126").Value = 39990.477
$ws.Range("N126").Value = -49870.477

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
# Row 31
$ws.Range("H31").Value = 1807.5834
$ws.Range("I31").Value = 1370.2222
$ws.Range("K31").Value = 1370.2222
$ws.Range("M31").Value = -1075.2222
# Row 34
$ws.Range("H34").Value = 1807.5834
$ws.Range("I34").Value = 1370.2222
$ws.Range("K34").Value = 1370.2222
$ws.Range("M34").Value = -1168.2222
# Row 99
$ws.Range("H99").Value = 2333.3333
$ws.Range("J99").Value = 2500
$ws.Range("L99").Value = 2500
$ws.Range("N99").Value = -5496
# Row 126
$ws.Range("H126").Value = 2333.3333
$ws.Range("J126").Value = 2500
$ws.Range("L126").Value = 7500
$ws.Range("N126").Value = -12440
# Row 132
$ws.Range("H132").Value = 1544.6428
$ws.Range("I132").Value = 863.14703
$ws.Range("K132").Value = 2589.44109
$ws.Range("M132").Value = -59.4410899999998

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 87.85714
$ws.Range("I2").Value = 106.5
$ws.Range("K2").Value = 639
$ws.Range("M2").Value = -526
# Row 17
$ws.Range("H17").Value = 166669150
$ws.Range("J17").Value = 2976
$ws.Range("L17").Value = 8928
$ws.Range("N17").Value = -9266
# Row 29
$ws.Range("H29").Value = 134.88889
$ws.Range("J29").Value = 150.33333
$ws.Range("L29").Value = 450.99999
$ws.Range("N29").Value = -1004.99999
# Row 34
$ws.Range("H34").Value = 717.375
$ws.Range("J34").Value = 1005.9
$ws.Range("L34").Value = 3017.7
$ws.Range("N34").Value = -3185.7
# Row 39
$ws.Range("H39").Value = 3542.5
$ws.Range("J39").Value = 3877.1428
$ws.Range("L39").Value = 11631.4284
$ws.Range("N39").Value = -12219.4284
# Row 55
$ws.Range("H55").Value = 3916.6667
$ws.Range("J55").Value = 3916.6667
$ws.Range("L55").Value = 11750.0001
$ws.Range("N55").Value = -12104.0001
# Row 108
$ws.Range("H108").Value = 3000
$ws.Range("I108").Value = 1000
$ws.Range("K108").Value = 3000
$ws.Range("M108").Value = -120
# Row 113
$ws.Range("H113").Value = 7530.467
$ws.Range("I113").Value = 33833.668
$ws.Range("K113").Value = 101501.004
$ws.Range("M113").Value = -99331.00399999999

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 19379.8
$ws.Range("I70").Value = 42000
$ws.Range("J70").Value = 4299.6665
$ws.Range("K70").Value = 42000
$ws.Range("L70").Value = 4299.6665
$ws.Range("M70").Value = -41730
$ws.Range("N70").Value = -4839.6665
# Row 73
$ws.Range("H73").Value = 19379.8
$ws.Range("I73").Value = 42000
$ws.Range("J73").Value = 4299.6665
$ws.Range("K73").Value = 42000
$ws.Range("L73").Value = 4299.6665
$ws.Range("M73").Value = -41064
$ws.Range("N73").Value = -6171.6665
# Row 102
$ws.Range("H102").Value = 3560.25
$ws.Range("I102").Value = 3711.7144
$ws.Range("K102").Value = 3711.7144
$ws.Range("M102").Value = -2089.7144
# Row 107
$ws.Range("H107").Value = 1433.3334
$ws.Range("J107").Value = 4000
$ws.Range("L107").Value = 4000
$ws.Range("N107").Value = -7840
# Row 123
$ws.Range("H123").Value = 14691
$ws.Range("J123").Value = 14691
$ws.Range("L123").Value = 14691
$ws.Range("N123").Value = -19591
# Row 126
$ws.Range("H126").Value = 2573080.2
$ws.Range("I126").Value = 11114173
$ws.Range("K126").Value = 33342519
$ws.Range("M126").Value = -33340049
# Row 132
$ws.Range("H132").Value = 714302.4
$ws.Range("I132").Value = 1069816.9
$ws.Range("K132").Value = 3209450.7
$ws.Range("M132").Value = -3206920.7
# Row 134
$ws.Range("H134").Value = 24460.889
$ws.Range("J134").Value = 24460.889
$ws.Range("L134").Value = 73382.667
$ws.Range("N134").Value = -78452.667

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2782.6316
$ws.Range("I7").Value = 2225.25
$ws.Range("K7").Value = 2225.25
$ws.Range("M7").Value = -2113.25
# Row 40
$ws.Range("H40").Value = 7323.6523
$ws.Range("I40").Value = 7734.4375
$ws.Range("K40").Value = 7734.4375
$ws.Range("M40").Value = -7598.4375
# Row 93
$ws.Range("H93").Value = 1391.1111
$ws.Range("I93").Value = 805.3333
$ws.Range("J93").Value = 2562.6667
$ws.Range("K93").Value = 805.3333
$ws.Range("L93").Value = 2562.6667
$ws.Range("M93").Value = 442.6667
$ws.Range("N93").Value = -5058.6667
# Row 126
$ws.Range("H126").Value = 2782.6316
$ws.Range("I126").Value = 2225.25
$ws.Range("K126").Value = 6675.75
$ws.Range("M126").Value = -4205.75
# Row 132
$ws.Range("H132").Value = 1274.0405
$ws.Range("I132").Value = 870.8302
$ws.Range("J132").Value = 2291.6667
$ws.Range("K132").Value = 2612.4906
$ws.Range("L132").Value = 6875.000100000001
$ws.Range("M132").Value = -82.49060000000009
$ws.Range("N132").Value = -11935.0001
# Row 136
$ws.Range("H136").Value = 1319.5375
$ws.Range("J136").Value = 3758.4
$ws.Range("L136").Value = 11275.2
$ws.Range("N136").Value = -16375.2

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 124
$ws.Range("H124").Value = 23747.5
$ws.Range("J124").Value = 23747.5
$ws.Range("L124").Value = 23747.5
$ws.Range("N124").Value = -33567.5
# Row 126
$ws.Range("H126").Value = 6785.6816
$ws.Range("I126").Value = 7543.25
$ws.Range("J126").Value = 4765.5
$ws.Range("K126").Value = 22629.75
$ws.Range("L126").Value = 14296.5
$ws.Range("M126").Value = -20159.75
$ws.Range("N126").Value = -19236.5
